$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RESOURCES")

# Insert a new column F ("costs_sell_kWh" - selling price) before the existing
# "reference" column (which shifts from F to G, carrying its data/styles along).
$ws.Columns("F:F").Insert()

# Header for the new column.
$ws.Range("F1").Value = "costs_sell_kWh"

# Selling-price values/formulas, mirroring the adjacent buying-cost column (E).
$ws.Range("F2").Value = 0
$ws.Range("F3").Formula = "=4.94/293"
$ws.Range("F4").Formula = "=0.2*0.75"
$ws.Range("F5").Value = 0.0001

# Header comment describing the new column.
$ws.Range("F1").AddComment('Selling price in US$(2015)/kWh(resource [thermal in case of fuels]).yr') | Out-Null

# Match the workbook's recorded selection.
$ws.Range("F1").Select() | Out-Null
